$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Product of Array Except Self" row, added below existing data (row 11)
$ws.Range("A11").Value = "Product of Array Except Self"
$ws.Range("B11").Value = "Return array of product of each element except current"
$ws.Range("C11").Value = "Create prefix and postfix product array. Calculate prefix and postfix products and return array by multiplying both"

# Match styling used by other rows (Name column uses the "Neutral" style,
# Description/Approach columns use the default body style) by copying the
# formats from the row above, reusing existing style entries.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection as recorded in the saved workbook
$ws.Range("C12").Select()
